$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing quarterly data (old D:K) to F:M
$ws.Columns("D:E").Insert()

# Copy the number formatting from the shifted columns (F:G, which carry the old D:E format)
# into the newly inserted D:E columns so the new cells match the original column styling
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarters (D = Q4 2018, E = Q3 2018) with reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 57000
$ws.Range("E8").Value = 55200
$ws.Range("D9").Value = 13600
$ws.Range("E9").Value = 13300
$ws.Range("D10").Value = 43400
$ws.Range("E10").Value = 41900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 200
$ws.Range("E14").Value = 100
$ws.Range("D15").Value = 20700
$ws.Range("E15").Value = 20100
$ws.Range("D17").Value = 40800
$ws.Range("E17").Value = 39800
$ws.Range("D18").Value = 16200
$ws.Range("E18").Value = 15400
$ws.Range("D20").Value = 5600
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 42500
$ws.Range("E21").Value = 35600
$ws.Range("D22").Value = 6700
$ws.Range("E22").Value = 6500
$ws.Range("D23").Value = 15200
$ws.Range("E23").Value = 9000
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 15200
$ws.Range("E26").Value = 9000
$ws.Range("D27").Value = 12400
$ws.Range("E27").Value = 6300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -5600
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 12400
$ws.Range("E33").Value = 6300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 12400
$ws.Range("E35").Value = 6300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 180600
$ws.Range("E41").Value = 183900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 28100
$ws.Range("E47").Value = 27100
$ws.Range("D48").Value = 2487300
$ws.Range("E48").Value = 2363000
$ws.Range("D49").Value = 60800
$ws.Range("E49").Value = 58600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1300
$ws.Range("E52").Value = 1500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2787700
$ws.Range("E54").Value = 2668800
$ws.Range("D57").Value = 21100
$ws.Range("E57").Value = 30400
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 39200
$ws.Range("E59").Value = 37100
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 757400
$ws.Range("E61").Value = 757200
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 911600
$ws.Range("E66").Value = 914300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 159100
$ws.Range("E70").Value = 159100
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -88300
$ws.Range("E72").Value = -85400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1717000
$ws.Range("E76").Value = 1595500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 12400
$ws.Range("E81").Value = 6300
$ws.Range("D83").Value = 20700
$ws.Range("E83").Value = 20100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 21800
$ws.Range("E89").Value = 35800
$ws.Range("D91").Value = -16700
$ws.Range("E91").Value = -14800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -138700
$ws.Range("E94").Value = -48300
$ws.Range("D96").Value = -400
$ws.Range("E96").Value = -400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 113600
$ws.Range("E100").Value = 33700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -3300
$ws.Range("E102").Value = 21200
